# Update the data-quality summary numbers on Sheet1.
# REF universe went from 10 to 9 ids, NEW universe from 8 to 7 ids, and the
# "duplicates (ref)" breakdown changed from 1 dup (id 5) to 2 dups (ids 1,5),
# which cascades into the BOTH/GAP/DUPS rows and their percentages.
#
# Note: the target values are strings (e.g. "9", "15.4%") stored as text in
# the sheet, not numbers. Typing a bare numeric-looking string into a cell
# makes Excel auto-convert it to a real number, so those are entered with a
# leading apostrophe to force a text entry (this mirrors typing '9 by hand
# in Excel); the apostrophe itself is not part of the stored value/text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'9"
$ws.Range("D3").Value = "'7"
$ws.Range("D5").Value = "'15"
$ws.Range("D6").Value = "'2"
$ws.Range("F6").Value = "'15.4%"
$ws.Range("D7").Value = "'0"
$ws.Range("E7").Value = "'9"
$ws.Range("F7").Value = "'0.0%"
$ws.Range("E8").Value = "'6"
$ws.Range("F8").Value = "'33.3%"
$ws.Range("D9").Value = "'4"
$ws.Range("E9").Value = "'6"
$ws.Range("F9").Value = "'66.7%"
$ws.Range("G9").Value = "2,4,7,8"
$ws.Range("E10").Value = "'9"
$ws.Range("F10").Value = "'33.3%"
$ws.Range("E11").Value = "'7"
$ws.Range("F11").Value = "'14.3%"
$ws.Range("E12").Value = "'6"
$ws.Range("F12").Value = "'33.3%"
$ws.Range("D13").Value = "'2"
$ws.Range("E13").Value = "'9"
$ws.Range("F13").Value = "'22.2%"
$ws.Range("G13").Value = "1,5"
$ws.Range("E14").Value = "'7"
$ws.Range("F14").Value = "'14.3%"
